$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A106").Value = 'Senior Golang Developer'
$ws.Range("B106").Value = 'https://www.dice.com/job-detail/7d79dfc9-53f9-454a-aad8-ae3dcfc8ae7d'
$ws.Range("C106").Value = 'Plano, Texas'
$ws.Range("D106").Value = 'Contract'
$ws.Range("E106").Value = '$50 - $60'
$ws.Range("F106").Value = 'Mindlance'

$ws.Range("A107").Value = 'Senior Fullstack Golang Developer || Onsite at Phoenix,AZ & Plano,TX & Charlotte, NC || W2 & C2C'
$ws.Range("B107").Value = 'https://www.dice.com/job-detail/d7727f7c-2b77-4125-b2e4-e7c47d30162a'
$ws.Range("C107").Value = 'Phoenix, Arizona'
$ws.Range("D107").Value = 'Contract, Third Party'
$ws.Range("E107").Value = '55 - 60'
$ws.Range("F107").Value = 'NasTech Global, Inc.'
